$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.887.59'

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +2.00%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.880.37'

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.50%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.19%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '333.31'

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +3.39%  '

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.16%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4734'

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +5.64%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3972'

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +3.81%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '48.22'

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.11%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08055'

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +2.45%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.028'

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.51%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '21.90'

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +2.48%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.914.94'

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +2.70%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.964'

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +1.80%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.202'

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.13%  '

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.15%  '

$ws.Range("B17").Value = 'ShibaInu'

$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001051'

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +2.07%  '

$ws.Range("B18").Value = 'Litecoin'

$ws.Range("C18").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '87.29'

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.53%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06623'

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.36'

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +2.18%  '

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.08%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '28.033.49'

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +2.52%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.518'

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.97%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.06'

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +2.79%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.316'

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +2.59%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.133.53'

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +2.32%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '157.69'

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +3.95%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.28'

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +4.64%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.108'

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +2.33%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.622'

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +1.72%  '

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +2.39%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9845'

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +5.31%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09573'

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +2.73%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.464'

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.61%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.617'

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.30%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.329'

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +1.60%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06125'

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +2.80%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02265'

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +1.96%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.231'

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +2.36%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.253'

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.31%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6034'

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +2.35%  '

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.18%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1904'

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +3.07%  '

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +1.03%  '

$ws.Range("B45").Value = 'WEMIXTOKEN'

$ws.Range("C45").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.273'

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.92%  '

$ws.Range("B46").Value = 'Decentraland'

$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5721'

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +1.51%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '12.26'

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.72%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.417'

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.68%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.950'

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +1.23%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06834'

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.41%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '113.82'

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +5.24%  '
